$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.896.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.46%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.827.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.63%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.55%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'310.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.02%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.006"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.54%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4574"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.83%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3695"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.07160"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.18%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.8741"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.07%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07769"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.77%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'19.63"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.13%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.845.01"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'5.324"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.39%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'6.394"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.35%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'87.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.29%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'1.008"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.000008726"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.42%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'1.006"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.55%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'26.941.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.37%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'14.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.03%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "  -2.17%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'2.065.84"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'10.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.58%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'1.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.43%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'151.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'18.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.10%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'1.964"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.99%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'113.86"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.84%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'4.918"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.96%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.08788"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.76%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'3.031"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.26%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.7496"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.07%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'4.487"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.134"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.11%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'2.557"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.60%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "  +1.22%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.01946"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.90%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05141"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.30%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.906"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.65%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'6.952"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.40%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.4976"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.39%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.1599"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.53%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'8.327"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.55%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "  -3.11%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'1.006"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.58%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "  -1.93%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'102.11"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.99%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.614"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.34%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.06105"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.80%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'64.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.48%  "
$ws.Range("E51").Style = "Normal"
